$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44365
$ws.Range("J2").Value = 80
$ws.Range("D3").Value = 45177
$ws.Range("J3").Value = 120
$ws.Range("K3").Value = 26000
$ws.Range("L3").Value = 26000
$ws.Range("M3").Value = 26000
$ws.Range("P3").Value = 1733
$ws.Range("D4").Value = 44838
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 22000
$ws.Range("L4").Value = 22000
$ws.Range("M4").Value = 22000
$ws.Range("P4").Value = 1467
$ws.Range("D5").Value = 45156
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 25000
$ws.Range("L5").Value = 25000
$ws.Range("M5").Value = 25000
$ws.Range("P5").Value = 1667
$ws.Range("D6").Value = 44400
$ws.Range("J6").Value = 80
$ws.Range("K6").Value = 25000
$ws.Range("L6").Value = 25000
$ws.Range("M6").Value = 25000
$ws.Range("P6").Value = 1667
$ws.Range("D7").Value = 45163
$ws.Range("J7").Value = 140
$ws.Range("K7").Value = 25000
$ws.Range("L7").Value = 25000
$ws.Range("M7").Value = 25000
$ws.Range("P7").Value = 1667
$ws.Range("D8").Value = 44775
$ws.Range("J8").Value = 120
$ws.Range("D9").Value = 44771
$ws.Range("J9").Value = 90
$ws.Range("K9").Value = 25000
$ws.Range("L9").Value = 25000
$ws.Range("M9").Value = 25000
$ws.Range("P9").Value = 1667
$ws.Range("D10").Value = 44407
$ws.Range("J10").Value = 90
$ws.Range("D11").Value = 44789
$ws.Range("J11").Value = 90
$ws.Range("K11").Value = 24000
$ws.Range("L11").Value = 24000
$ws.Range("M11").Value = 24000
$ws.Range("P11").Value = 1600
$ws.Range("D12").Value = 45146
$ws.Range("J12").Value = 140
$ws.Range("K12").Value = 26000
$ws.Range("L12").Value = 26000
$ws.Range("M12").Value = 26000
$ws.Range("P12").Value = 1733
$ws.Range("D13").Value = 44799
$ws.Range("J13").Value = 80
$ws.Range("K13").Value = 23000
$ws.Range("L13").Value = 23000
$ws.Range("M13").Value = 23000
$ws.Range("P13").Value = 1533
$ws.Range("D14").Value = 44831
$ws.Range("D15").Value = 45149
$ws.Range("D16").Value = 44806
$ws.Range("J16").Value = 70
$ws.Range("K16").Value = 23000
$ws.Range("L16").Value = 23000
$ws.Range("M16").Value = 23000
$ws.Range("P16").Value = 1533
$ws.Range("D17").Value = 44757
$ws.Range("J17").Value = 80
$ws.Range("K17").Value = 25000
$ws.Range("L17").Value = 25000
$ws.Range("M17").Value = 25000
$ws.Range("P17").Value = 1667
$ws.Range("D18").Value = 44827
$ws.Range("K18").Value = 22000
$ws.Range("L18").Value = 22000
$ws.Range("M18").Value = 22000
$ws.Range("P18").Value = 1467
$ws.Range("D19").Value = 44418
$ws.Range("K19").Value = 25000
$ws.Range("L19").Value = 25000
$ws.Range("M19").Value = 25000
$ws.Range("P19").Value = 1667
$ws.Range("D20").Value = 44803
$ws.Range("J20").Value = 90
$ws.Range("D21").Value = 44761
$ws.Range("J21").Value = 100
$ws.Range("K21").Value = 23000
$ws.Range("L21").Value = 25000
$ws.Range("M21").Value = 24000
$ws.Range("P21").Value = 1600
$ws.Range("D22").Value = 44792
$ws.Range("J22").Value = 120
$ws.Range("K22").Value = 24000
$ws.Range("L22").Value = 24000
$ws.Range("M22").Value = 24000
$ws.Range("P22").Value = 1600
$ws.Range("D23").Value = 44781
$ws.Range("J23").Value = 70
$ws.Range("K23").Value = 24000
$ws.Range("L23").Value = 24000
$ws.Range("M23").Value = 24000
$ws.Range("P23").Value = 1600
$ws.Range("D24").Value = 44817
$ws.Range("K24").Value = 23000
$ws.Range("L24").Value = 23000
$ws.Range("M24").Value = 23000
$ws.Range("P24").Value = 1533
$ws.Range("D25").Value = 44819
$ws.Range("J25").Value = 70
$ws.Range("K25").Value = 22000
$ws.Range("L25").Value = 22000
$ws.Range("M25").Value = 22000
$ws.Range("P25").Value = 1467
$ws.Range("D26").Value = 44810
$ws.Range("J26").Value = 110
$ws.Range("K26").Value = 22000
$ws.Range("L26").Value = 22000
$ws.Range("M26").Value = 22000
$ws.Range("P26").Value = 1467
$ws.Range("D27").Value = 44740
$ws.Range("D28").Value = 45167
$ws.Range("J28").Value = 120
$ws.Range("D29").Value = 44782
$ws.Range("J29").Value = 120
$ws.Range("K29").Value = 24000
$ws.Range("L29").Value = 24000
$ws.Range("M29").Value = 24000
$ws.Range("P29").Value = 1600
$ws.Range("D30").Value = 44754
$ws.Range("J30").Value = 90
$ws.Range("D31").Value = 44764
$ws.Range("K31").Value = 24000
$ws.Range("L31").Value = 24000
$ws.Range("M31").Value = 24000
$ws.Range("P31").Value = 1600
$ws.Range("D32").Value = 44750
$ws.Range("K32").Value = 25000
$ws.Range("L32").Value = 25000
$ws.Range("M32").Value = 25000
$ws.Range("P32").Value = 1667
$ws.Range("D33").Value = 44778
$ws.Range("J33").Value = 120
$ws.Range("K33").Value = 24000
$ws.Range("L33").Value = 24000
